$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "Matiere1"
$ws.Range("H1").Value = "Matiere2"

$ws.Range("H1").Select()
